# Grade.xlsx SPRING 2024 update
# - FML!B5 (Assignment 3 score) set to 15
# - AiP!B3 (Quiz 1 score) updated to 3.562
# - AiP!B4 (Quiz 2 score) updated to 3.635 (was a formula =C4, now a literal value)
# - AiP!B5 (Quiz 3 score) set to 3.197 (previously empty)
# - Active worksheet moves from AiP to FML, with B6 selected on both sheets

$wb = $excel.ActiveWorkbook

$wsFML = $wb.Worksheets.Item("FML")
$wsAiP = $wb.Worksheets.Item("AiP")

# --- Data edits -----------------------------------------------------------

# FML: Assignment 3 score entered
$wsFML.Range("B5").Value = 15

# AiP: quiz scores entered / corrected
$wsAiP.Range("B3").Value = 3.562
$wsAiP.Range("B4").Value = 3.635
$wsAiP.Range("B5").Value = 3.197

# --- Selection / active sheet ---------------------------------------------
# Match final selection state: AiP loses tabSelected, FML gains it; both
# sheets end up with B6 as the active cell / selection.

$wsAiP.Activate()
$wsAiP.Range("B6").Select()

$wsFML.Activate()
$wsFML.Range("B6").Select()
